# "updated script 2nd day"
# Replace the StartDate/StartTime/EndDate/EndTime/Location columns (C..G)
# with a smaller Field/Age table, and add a sample data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old headers/data that used to live in columns C..G
$ws.Range("C1:G2").ClearContents()

# New headers in row 1
$ws.Range("C1").Value = "Field"
$ws.Range("D1").Value = "Age"

# Row 2 keeps the existing Hello / Hellodesc values (A2/B2) and gains
# a sample row for the new Field/Age columns
$ws.Range("C2").Value = "TestSample"
$ws.Range("D2").Value = 30

# Give the new "Field" column a dedicated width, like "Description" has
$ws.Columns.Item(3).ColumnWidth = 10.33

# Selection ends up on the newly entered Age value
$ws.Range("D2").Select()
